$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.764.20"
$ws.Range("D3").Value = "1.547.97"
$ws.Range("E3").Value = "  -1.84%  "
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "204.41"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -1.78%  "
$ws.Range("E6").Value = "  -1.76%  "
$ws.Range("E7").Value = "  +0.04%  "
$ws.Range("E8").Value = "  -1.23%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "21.36"
$ws.Range("D9").ClearFormats()
$ws.Range("E10").Value = "  -1.82%  "
$ws.Range("E11").Value = "  -0.88%  "
$ws.Range("D12").Value = "1.768.99"
$ws.Range("E12").Value = "  -1.82%  "
$ws.Range("D13").Value = "1.545.91"
$ws.Range("E13").Value = "  -2.05%  "
$ws.Range("E14").Value = "  -2.80%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.510"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  -1.87%  "
$ws.Range("D16").Value = "26.779.00"
$ws.Range("E16").Value = "  -1.93%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "60.86"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  -2.68%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "213.55"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  -0.84%  "
$ws.Range("D19").Value = "0.0₃0684"
$ws.Range("E19").Value = "  -0.53%  "
$ws.Range("E20").Value = "  -1.63%  "
$ws.Range("E21").Value = "  +0.05%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.07"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  -1.77%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "9.01"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  -4.42%  "
$ws.Range("E24").Value = "  -1.09%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "152.68"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  +0.50%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "6.52"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  -2.57%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "14.84"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  -0.84%  "
$ws.Range("E28").Value = "  +0.04%  "
$ws.Range("E29").Value = "  -2.15%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.0460"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  -0.99%  "
$ws.Range("E31").Value = "  -3.70%  "
$ws.Range("E32").Value = "  -0.35%  "
$ws.Range("D33").Value = "1.351.37"
$ws.Range("E33").Value = "  -4.25%  "
$ws.Range("E34").Value = "  -1.02%  "
$ws.Range("E35").Value = "  -3.75%  "
$ws.Range("E36").Value = "  -0.81%  "
$ws.Range("E37").Value = "  -2.49%  "
$ws.Range("E38").Value = "  -2.14%  "
$ws.Range("E39").Value = "  +0.51%  "
$ws.Range("E40").Value = "  -2.51%  "
$ws.Range("E41").Value = "  +0.02%  "
$ws.Range("B42").Value = "WEMIXToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.989"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  -1.13%  "
$ws.Range("B43").Value = "FraxShare"
$ws.Range("C43").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.54"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  +3.78%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.19"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  +0.18%  "
$ws.Range("E45").Value = "  -3.10%  "
$ws.Range("E46").Value = "  -1.77%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.28"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  -3.35%  "
$ws.Range("D48").Value = "1.683.04"
$ws.Range("E48").Value = "  -1.82%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "85.85"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  -0.36%  "
$ws.Range("E50").Value = "  +2.33%  "
$ws.Range("E51").Value = "  -1.38%  "
